# Update "想去人数" (interest count) figures in the 展览 and 全部类型 sheets.
# F2: 607 -> 608
# F3: 561 -> 562
# F6: 94  -> 96
# F10: 4887 -> 4892
# F11: 4604 -> 4607

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 608
    $ws.Range("F3").Value = 562
    $ws.Range("F6").Value = 96
    $ws.Range("F10").Value = 4892
    $ws.Range("F11").Value = 4607
}
